$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-30 Wednesday" "2025-07-31 Thursday"

Replace-Text "360÷3=" "939÷9="
Replace-Text "868÷5=" "689÷5="
Replace-Text "234÷9=" "561÷9="
Replace-Text "315÷2=" "821÷2="
Replace-Text "334÷4=" "556÷8="

Replace-Text "559÷3=" "184÷6="
Replace-Text "225÷2=" "965÷6="
Replace-Text "981÷4=" "998÷6="
Replace-Text "746÷5=" "269÷2="
Replace-Text "600÷6=" "175÷4="

Replace-Text "484÷3=" "377÷6="
Replace-Text "762÷3=" "973÷4="
Replace-Text "996÷3=" "170÷9="
Replace-Text "450÷3=" "514÷8="
Replace-Text "894÷9=" "476÷6="

Replace-Text "779÷7=" "571÷4="
Replace-Text "405÷7=" "130÷8="
Replace-Text "515÷7=" "839÷9="
Replace-Text "100÷7=" "718÷3="
Replace-Text "638÷2=" "403÷7="

Replace-Text "242÷7=" "824÷9="
Replace-Text "926÷6=" "583÷2="
Replace-Text "821÷4=" "714÷2="
Replace-Text "273÷2=" "913÷3="
Replace-Text "314÷9=" "215÷8="
